# PP & test plan update
#
# Summary of the change (see commit diff):
#  - Row 9  (TestNormalGRP)      : Status "In progress" -> "Ran - Failed"; add Comments "Printing error"
#  - Row 10 (TestNormalPARAGRP)  : Status "In progress" -> "Ran - Failed"; add Comments "Testing specific
#                                   error, not 'real' error" (wrapped, like the other Comments cells)
#  - Remove the two trailing blank rows (13 and 14)
#  - Move the active selection to H13 (where the next entry would go)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 first, so the new shared string for its comment is inserted
#     before the one used by row 9 (matches the order new <si> entries show
#     up in xl/sharedStrings.xml). ---
$ws.Range("G10").Value = "Ran - Failed"
$ws.Range("H10").Value = "Testing specific error, not 'real' error"
$ws.Range("H10").WrapText = $true
$ws.Rows(10).RowHeight = 28.8

$ws.Range("G9").Value = "Ran - Failed"
$ws.Range("H9").Value = "Printing error"

# Remove the now-obsolete blank rows at the bottom of the sheet.
$ws.Rows("13:14").Delete()

# Match the saved selection/active cell.
$ws.Range("H13").Select()
